$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (price/volume columns). Cells store
# numeric-looking text, so a leading apostrophe forces text entry;
# Style is reset to Normal afterwards so no formatting/quote-prefix
# artifact is left behind on the cell.
$ws.Range("D2").Value = "'301.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.18%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'2.40%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.986"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.19%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07740"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.39%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.073"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-4.91%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.907"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.53%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9226"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.62%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.09699"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'5.39%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1860"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.96%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08600"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.77%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03513"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.84%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09930"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.41%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001464"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-3.08%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005616"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.73%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.462"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.49%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.029"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.17%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.417"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'10.71%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3407"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.61%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1342"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.23%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.756"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.61%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2196"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.32%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04584"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.23%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'13.29%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001228"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.92%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001397"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'6.64%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E39").Value = "'2.36%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04638"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.26%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007430"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-6.56%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1388"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.16%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007704"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.11%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002154"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-3.40%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01031"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'14.78%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006181"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.01%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.93%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0005790"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-0.19%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'35.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'570.90%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-26.36%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.00002095"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.93%"
$ws.Range("E51").Style = "Normal"
